$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: update the Number of Test Cases value
$ws.Range("F7").Value = 43

# Row 8: new data row - (TS_005) Header_Create Group
$ws.Range("B8").Value = "(TS_005)" + [char]10 + "Header_Create Group"
$ws.Range("C8").Value = "FRS"
$ws.Range("D8").Value = "Validte the functionality of the Dropdown Header > Create Group."
$ws.Range("F8").Value = 13
$ws.Rows.Item(8).RowHeight = 30

# Row 9: new data row - (TS_006) Header_Fairness Calculator
$ws.Range("B9").Value = "(TS_006)" + [char]10 + "Header_Fairness Calculator"
$ws.Range("C9").Value = "FRS"
$ws.Range("D9").Value = "Validte the functionality of the Dropdown Header > Fairness Calculator."
$ws.Rows.Item(9).RowHeight = 45

# Row 10: move the SUM formula down to F10 (pick up F9's cell formatting, since
# F10 had no cell/style of its own before)
$ws.Range("F9").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F10").Formula = "=SUM(F4:F9)"

# Update selection to F11
$ws.Range("F11").Select()
